$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the plan-code values used by the rule rows so that the campaign
# now also applies to the RBWLI3A plan.
$ws.Range("D10:D27").Value = "RBWLI3, RBWLI3A"

# Update the rule-table "PLAN CODE" condition to use a contains() check
# instead of an exact equality check.
$ws.Range("D8").Value = "planCode.contains($1)"

# Column D now needs to be a bit wider to comfortably fit the longer text,
# so split it out from the shared B:D column-width group.
$ws.Columns("D").ColumnWidth = 20.1

# Update the saved cursor/selection position.
[void]$ws.Range("F35").Select()
